# Remove the duplicated English-language "Information on the processing of
# personal data" section (added after the Italian "Informativa" text by a
# manual page-break paragraph) so only the Italian version remains.
#
# Structure before the edit (within the main document body):
#   ... Italian section ...
#   <empty paragraph containing only a manual page break>
#   "Information on the processing of personal data" (English title) ...
#   ... the rest of the English translation, section 1..9 ...
#   "...or to take appropriate legal action (art. 79 of the GDPR)."  <- last paragraph
#   <w:sectPr> (end of body)
#
# Structure after the edit:
#   ... Italian section ...
#   <empty paragraph, no page break, directly followed by sectPr>

$d = $word.ActiveDocument

# Locate (by content, not by hard-coded index) the paragraph that starts the
# English section and the paragraph that ends it, using the global
# Paragraphs collection (whose Range.Start/End are reliable, unlike ranges
# re-derived from a sub-range).
$all = $d.Paragraphs
$count = $all.Count

$englishTitleIndex = -1
$englishLastIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $para = $all.Item($i)
    $text = $para.Range.Text
    if ($text -match "^Information on the processing of personal data") {
        $englishTitleIndex = $i
    }
    if ($text -match "as provided for by art\. 77 of the GDPR itself") {
        $englishLastIndex = $i
    }
}

if ($englishTitleIndex -eq -1 -or $englishLastIndex -eq -1) {
    throw "Could not locate the English privacy-notice section to remove."
}

# The paragraph right before the English title is the one holding the manual
# page break that separates the Italian and English copies.
$breakIndex = $englishTitleIndex - 1
$breakParagraph = $d.Paragraphs.Item($breakIndex)
$breakParaStart = $breakParagraph.Range.Start

# Delete every paragraph of the English section, from its first paragraph
# through (and including) its very last paragraph / paragraph mark. This
# merges what used to be the start of the English title paragraph into the
# break paragraph, leaving the break paragraph as the final (now empty)
# paragraph of the body, immediately followed by the section properties.
$englishStart = $d.Paragraphs.Item($englishTitleIndex).Range.Start
$englishEnd = $d.Paragraphs.Item($englishLastIndex).Range.End
$englishRange = $d.Range($englishStart, $englishEnd)
$englishRange.Delete()

# Remove the manual page break character left in the (previously) break-only
# paragraph, turning it into a plain empty paragraph while preserving its own
# paragraph formatting (spacing/jc/rPr already on that paragraph mark).
$breakCharRange = $d.Range($breakParaStart, $breakParaStart + 1)
if ($breakCharRange.Text -eq [char]12) {
    $breakCharRange.Delete()
}
